$d = $word.ActiveDocument

$replacements = @(
    @('999×9=8991', '472×4=1888'),
    @('906×3=2718', '682×8=5456'),
    @('405×9=3645', '899×7=6293'),
    @('743×2=1486', '892×6=5352'),
    @('724×4=2896', '678×5=3390'),
    @('254×4=1016', '699×9=6291'),
    @('492×8=3936', '383×6=2298'),
    @('810×8=6480', '939×2=1878'),
    @('159×8=1272', '740×7=5180'),
    @('317×6=1902', '125×8=1000'),
    @('558×2=1116', '434×4=1736'),
    @('947×4=3788', '546×4=2184'),
    @('998×4=3992', '267×7=1869'),
    @('373×2=746', '864×5=4320'),
    @('184×8=1472', '582×5=2910'),
    @('474×7=3318', '294×7=2058'),
    @('762×4=3048', '300×4=1200'),
    @('679×3=2037', '245×5=1225'),
    @('488×9=4392', '378×4=1512'),
    @('568×8=4544', '982×2=1964'),
    @('499×8=3992', '741×9=6669'),
    @('352×7=2464', '612×6=3672'),
    @('793×7=5551', '873×9=7857'),
    @('787×2=1574', '410×9=3690'),
    @('457×5=2285', '599×2=1198'),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

$d.Save()
